# Allow GC's to hide medicare fields
#
# The "RequireMedical" field is being turned into a "MedicareOption" field
# that (in addition to TRUE/FALSE) can also be set to "Hide" so group
# coordinators can hide medicare fields on the signup form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column Q header from "RequireMedical" to "MedicareOption".
$ws.Range("Q1").Value = "MedicareOption"

# Update the sample/default row value for that field from "FALSE" to "Hide".
$ws.Range("Q2").Value = "Hide"

# Match the author's recorded selection/view state at the time of the edit.
[void]$ws.Range("Q3").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1

Write-Output "MedicareOption/Hide update applied"
